$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clusters = @(
  '3323 Villa Maria Catholic Homes St Bernadette''s Aged Care Sunshine North',
  '3364 Assisi Centre Aged Care Rosanna',
  '3376 Royal Freemasons Coppin Centre Melbourne',
  '3622 Olivet Care Aged Care Services Ringwood',
  '3825 TLC Forest Lodge Residential Aged Care Frankston North',
  '3961 Water Gardens Aged Care Sydenham Tier 1A',
  '4167 Royal Freemasons Centennial Lodge Wantirna South',
  '4282 Villa Maria Catholic Homes (VMCH) Wantirna At-Home Aged Care',
  '45034 River Gum Primary School Hampton Park',
  '45573 Narre Warren South P-12 College Narre Warren South',
  '45695 Sacred Heart Primary School Yarrawonga',
  '50516 Ilim College Glenroy Campus Hadfield',
  '50567 Alamanda K9 College Point Cook',
  '52912 Edgars Creek Primary School Wollert',
  '52985 Minaret College Springvale',
  'AG Industries Pty Ltd Factory Thomastown',
  'Adass Israel School Elsternwick',
  'Antonine College Cedar Campus Coburg',
  'Bacchus Marsh Childcare and Kindergarten Centre Bacchus Marsh',
  'Baden Powell College Tarneit',
  'Collingwood College Abbotsford',
  'Covenant College Bell Post Hill',
  'Dandenong South Primary School Dandenong',
  'Darul Ulum College of Victoria Fawkner October',
  'Derrimut Primary School Derrimut',
  'Devon Meadows Primary School Devon Meadows',
  'Exford Primary School Exford',
  'Flemington Racecourse Flemington',
  'Gilly''s Early Learning Centre Balaclava',
  'Guardian Childcare & Education Moorabbin',
  'Hazel Glen College Doreen',
  'Hazelwood North Primary School Hazelwood North',
  'Ilim College Dallas Main Campus Dallas Oct',
  'Islamic College of Melbourne Tarneit Oct Nov',
  'Lyndhurst Primary School Lyndhurst',
  'Master Poultry Group West Footscray',
  'Middle Park Primary School Middle Park',
  'Minaret College Officer Campus Officer',
  'Morwell Park Primary School Morwell',
  'Nido Early School Woodend',
  'Nio Early Learning Adventures Preston',
  'Oakleigh South Primary School Oakleigh South',
  'Pentland Primary School Darley',
  'Rutherglen Motor Inn and Walkabout Motel Rutherglen',
  'Sirius College Ibrahim Dellal Campus Sunshine',
  'Sirius College Shepparton Campus Shepparton',
  'Smartie Pants Early Learning and Development Diamond Creek',
  'Social Gathering Woodvale 30 Oct',
  'Society Restaurant Melbourne',
  'St Ambrose Parish Primary School Woodend',
  'St Clare''s Primary School Officer',
  'St Georges Road Primary School Shepparton',
  'St Joseph''s School Quarry Hill',
  'St Louis de Montfort''s School Aspendale',
  'St Paul''s Primary School Sunshine West',
  'Stevensville Primary School St Albans',
  'Stockdale Road Primary School Traralgon',
  'Supreme Caravans Manufacturing Campbellfield',
  'Templestowe Park Primary School Templestowe',
  'The Lake Primary School Cabarita',
  'Top Yard Rooftop Melbourne',
  'Truganina P-9 College Truganina',
  'Tucker Road Bentleigh Primary School Bentleigh',
  'Warragul Regional College Warragul',
  'Wodonga Primary School Wodonga',
  'Wyndham Christian College Wyndham Vale',
  'Yeshivah College St Kilda East'
)

$counts = @(
  11,
  28,
  30,
  11,
  15,
  20,
  23,
  11,
  10,
  17,
  28,
  11,
  13,
  10,
  12,
  13,
  12,
  14,
  25,
  10,
  10,
  26,
  15,
  11,
  11,
  12,
  16,
  13,
  11,
  13,
  14,
  23,
  12,
  48,
  14,
  12,
  11,
  29,
  59,
  14,
  15,
  11,
  12,
  22,
  14,
  22,
  17,
  10,
  27,
  12,
  11,
  15,
  32,
  10,
  19,
  11,
  30,
  50,
  33,
  23,
  14,
  12,
  11,
  19,
  12,
  11,
  24
)

for ($i = 0; $i -lt $clusters.Count; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $clusters[$i]
  $ws.Cells.Item($row, 2).Value = $counts[$i]
}
